$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0935887046539458, 0.9633568397822009, 4.680461315292987, 2.163437384185867, 2.191226136823886, 37)
    3  = @(0.1801871629257656, 0.926508908720455,  4.517135635143278, 2.125355413841007, 2.147743373686913, 36)
    4  = @(0.2352243557214495, 0.873084372529863,  4.374556468617144, 2.091544039368319, 2.108616140758171, 35)
    5  = @(0.1761038177183508, 0.918841548698121,  4.561172687284334, 2.135690213323162, 2.160425354518692, 34)
    6  = @(0.2382491768315642, 0.9330133540749047, 4.689726382363714, 2.165577609406718, 2.185805075505592, 33)
    7  = @(0.2167218967521913, 0.9364529417776732, 4.7812919825397,   2.186616560474127, 2.210665899271447, 32)
    8  = @(0.2578528310586296, 0.985525993686524,  4.976590512716229, 2.230827315754456, 2.252503609595107, 31)
    9  = @(0.2218826356104003, 0.9820920796201497, 5.119605115971289, 2.262654440247403, 2.290243145950139, 30)
    10 = @(0.2514602257437938, 1.021614200504279,  5.321885787980542, 2.306921279103503, 2.33376575473015,  29)
    11 = @(0.2366041340794119, 1.012053285505367,  5.483661848752136, 2.341721983659063, 2.372489310826974, 28)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
